# Apply the "Add files via upload" update to the Case Locations / Public
# Exposure Sites sheet:
#   - Insert a new data row at row 6 (Broadmeadows / Sacca's Fruit World),
#     pushing the former rows 6-38 down to 7-39.
#   - Normalize a handful of date/text typos that were corrected at the
#     same time.
#   - Rewrite the Queen Victoria Market note (now at D16) with the more
#     detailed "sheds A and B" wording.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 6 (Woolworths Broadmeadows),
# shifting it and everything below down by one row.
$ws.Rows("6:6").Insert()

# Fix the date typo on the row directly above the insertion point.
$ws.Range("C5").Value = "1:25pm - 1:59pm  9/2/2021"

# Populate the newly inserted row 6 with the new venue entry.
$ws.Range("A6").Value = "Broadmeadows"
$ws.Range("B6").Value = "Sacca's Fruit World  Broadmeadows Central  Broadmeadows VIC 3047"
$ws.Range("C6").Value = "12:30pm - 1:00pm  9/2/2021"
$ws.Range("D6").Value = "Case visited venue"

# Minor corrections on the rows that shifted down from their old
# positions (old row -> new row = old + 1).
$ws.Range("C7").Value = "12:15pm - 12:30pm 9/2/2021"
$ws.Range("C9").Value = "7:14pm  11:30pm  6/2/2021"
$ws.Range("C12").Value = "1:35pm  2:17pm  9/2/2021"
$ws.Range("B13").Value = "Melbourne Golf Academy  385 Centre Dandenong Rd  Heatherton VIC 3202"
$ws.Range("C14").Value = "6:45am - 7:30am  8/2/21"
$ws.Range("D16").Value = "Case attended sheds A and B (also known as section 2) - Fruit and Vegetables, and used female toilets adjacent to shed A.  See a map of the Queen Victoria Market (PDF)"
